$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("operador")
$ws2 = $wb.Worksheets.Item("destinos")

# -----------------------------------------------------------------
# Pre-format the new rows (format-only copy does not create shared
# strings, so this can safely happen before any value is written).
# Row 146 on "operador" is intentionally formatted (but not filled
# in) later, after the AutoFilter range is set - see note below.
# -----------------------------------------------------------------
$ws1.Range("A142:G142").Copy()
$ws1.Range("A143:G145").PasteSpecial(-4122) # xlPasteFormats

$ws2.Range("A273:I273").Copy()
$ws2.Range("A274:I275").PasteSpecial(-4122) # xlPasteFormats

# -----------------------------------------------------------------
# Fill in the new rows. The order in which new (never seen before)
# string values are written determines their position in the shared
# string table, so the writes below follow the same interleaved
# sheet1 / sheet2 sequence as the original edit.
# -----------------------------------------------------------------

# operador!143 - HiFly
$ws1.Cells.Item(143,1).Value = "3L"
$ws1.Cells.Item(143,2).Value = "HiFly"
$ws1.Cells.Item(143,3).Value = "HFM"
$ws1.Cells.Item(143,4).Value = "Hi Fly Malta"
$ws1.Cells.Item(143,5).Value = "Operator"
$ws1.Cells.Item(143,7).Value = "Sem Aliança"

# destinos!274 - Beja (Portugal)
$ws2.Cells.Item(274,1).Value = "BYJ"
$ws2.Cells.Item(274,2).Value = "PT"
$ws2.Cells.Item(274,3).Value = "Beja"
$ws2.Cells.Item(274,4).Value = "Portugal"
$ws2.Cells.Item(274,5).Value = "Europa"
$ws2.Cells.Item(274,6).ClearContents()
$ws2.Cells.Item(274,7).ClearContents()
$ws2.Cells.Item(274,8).Value = "Inter"
$ws2.Cells.Item(274,9).Value = "LPBJ"

# operador!144 - Martinair
$ws1.Cells.Item(144,1).Value = "MP"
$ws1.Cells.Item(144,2).Value = "Martinair"
$ws1.Cells.Item(144,3).Value = "MPH"
$ws1.Cells.Item(144,4).Value = "Martinair Holland NV"
$ws1.Cells.Item(144,5).Value = "Operator"
$ws1.Cells.Item(144,7).Value = "Sem Aliança"

# operador!145 - JetSmart (PER)
$ws1.Cells.Item(145,1).Value = "JZ"
$ws1.Cells.Item(145,2).Value = "JetSmart (PER)"
$ws1.Cells.Item(145,3).Value = "JAP"
$ws1.Cells.Item(145,4).Value = "JetSmart Airlines Perú S.A.C"
$ws1.Cells.Item(145,5).Value = "Operator"
$ws1.Cells.Item(145,7).Value = "Sem Aliança"

# destinos!275 - Araxá (Brasil)
$ws2.Cells.Item(275,1).Value = "AAX"
$ws2.Cells.Item(275,2).Value = "BR"
$ws2.Cells.Item(275,3).Value = "Araxá"
$ws2.Cells.Item(275,4).Value = "Brasil"
$ws2.Cells.Item(275,5).Value = "América do Sul"
$ws2.Cells.Item(275,6).Value = "Sudeste"
$ws2.Cells.Item(275,7).Value = "Minas Gerais"
$ws2.Cells.Item(275,8).Value = "Dom"
$ws2.Cells.Item(275,9).Value = "SBAX"

# -----------------------------------------------------------------
# Extend the AutoFilter range on "operador" to cover rows up to 145
# (turn it off first so the call doesn't just toggle the existing
# filter off), and keep the hidden _FilterDatabase name in sync.
# This must happen BEFORE row 146 is added below: AutoFilter always
# snaps to the full contiguous used range, so if row 146 already had
# data the filter would stretch to G146 instead of stopping at G145
# (matching the source file, where the filter was not extended to
# include the very last added row).
# "destinos" autoFilter / _FilterDatabase intentionally stay as-is.
# -----------------------------------------------------------------
$ws1.AutoFilterMode = $false
$ws1.Range("A1:G145").AutoFilter() | Out-Null

foreach ($n in $wb.Names) {
    if ($n.Name -eq "operador!_FilterDatabase") {
        $n.RefersTo = "=operador!`$A`$1:`$G`$145"
    }
}

# operador!146 - Anivia (added after the AutoFilter range is fixed)
$ws1.Range("A142:G142").Copy()
$ws1.Range("A146:G146").PasteSpecial(-4122) # xlPasteFormats
$ws1.Cells.Item(146,1).Value = "TOT"
$ws1.Cells.Item(146,2).Value = "Anivia"
$ws1.Cells.Item(146,3).Value = "TOT"
$ws1.Cells.Item(146,4).Value = "Anivia Serviços Aéreos LTDA"
$ws1.Cells.Item(146,5).Value = "Operator"
$ws1.Cells.Item(146,7).Value = "Sem Aliança"

# -----------------------------------------------------------------
# Leave the selection on each sheet at the first empty row below the
# newly added data, then return focus to "operador".
# -----------------------------------------------------------------
$ws2.Cells.Item(276,1).Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Cells.Item(147,1).Select() | Out-Null

Write-Host "done"
